$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "my_tasks" row (row 8): translation text changed from
# "My Tasks" / "Nhiệm vụ của tôi" to a survey-taking prompt.
$ws.Range("B8").Value = "Thực hiện khảo sát"
$ws.Range("C8").Value = "Take a pollution survey"

# Update the selected cell/range to match the saved view state.
$ws.Range("C19").Select()
